$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = 1593
$ws.Range("E20").Value = 53358905

$ws.Range("C26").Value = 33938
$ws.Range("E26").Value = 203949247

$ws.Range("C44").Value = 10556
$ws.Range("E44").Value = 42605161

$ws.Range("C74").Value = 951
$ws.Range("E74").Value = 4274646

$ws.Range("C77").Value = 4544
$ws.Range("E77").Value = 8328779

$ws.Range("C92").Value = 409150
$ws.Range("E92").Value = 1595480068

$ws.Range("C94").Value = 94199
$ws.Range("E94").Value = 917978553

$ws.Range("C95").Value = 50772
$ws.Range("E95").Value = 932793689

$ws.Range("C96").Value = 17294
$ws.Range("E96").Value = 794330219

$ws.Range("C97").Value = 2157
$ws.Range("E97").Value = 214111625

$ws.Range("C110").Value = 397
$ws.Range("E110").Value = 16698678

$ws.Range("C141").Value = 80474
$ws.Range("E141").Value = 280718664

$ws.Range("C174").Value = 226092
$ws.Range("E174").Value = 900594909
